# Add files via upload
# - Rename the existing sheet "データ1" -> "data"
# - Add a new sheet "setting" after "data" containing a short list of
#   time-of-day values (formatted as h:mm) in A1:A6
# - Keep "data" as the active sheet

$wb = $excel.ActiveWorkbook

# Rename the original (and only) worksheet.
$dataSheet = $wb.Worksheets.Item(1)
$dataSheet.Name = "data"

# Insert the new "setting" worksheet right after "data".
$settingSheet = $wb.Worksheets.Add($null, $dataSheet)
$settingSheet.Name = "setting"

# Populate the time values (stored as Excel day-fraction serials) and
# apply a time number format, matching the source workbook.
$settingSheet.Range("A1:A6").NumberFormat = "h:mm"
$settingSheet.Range("A1").Value = 0.35416666666666669
$settingSheet.Range("A2").Value = 0.5
$settingSheet.Range("A3").Value = 0.53472222222222221
$settingSheet.Range("A4").Value = 0.64583333333333337
$settingSheet.Range("A5").Value = 0.65277777777777779
$settingSheet.Range("A6").Value = 0.72222222222222221

# Leave "data" as the active/selected sheet.
$dataSheet.Activate()
